# A new daily price entry was added at the top of the data block (row 24),
# pushing every existing row from 24..58 down by one (to 25..59).
# The sheet's used range therefore grows from A1:T58 to A1:T59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 24; this shifts rows 24-58
# down to 25-59 and keeps everything else (formatting, other rows) intact.
$ws.Rows(24).Insert()

# Populate the newly inserted row 24 with the new record's values.
$ws.Cells.Item(24, 1).Value  = 10
$ws.Cells.Item(24, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(24, 3).Value  = 'La Araucanía'
$ws.Cells.Item(24, 4).Value2 = 44705
$ws.Cells.Item(24, 5).Value  = 9
$ws.Cells.Item(24, 6).Value  = 'Fruta'
$ws.Cells.Item(24, 7).Value  = 100107
$ws.Cells.Item(24, 8).Value  = 'Otros'
$ws.Cells.Item(24, 9).Value  = 100107011
$ws.Cells.Item(24, 10).Value = 'Tuna'
$ws.Cells.Item(24, 11).Value = 'Sin especificar'
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 25
$ws.Cells.Item(24, 14).Value = 20000
$ws.Cells.Item(24, 15).Value = 20000
$ws.Cells.Item(24, 16).Value = 20000
$ws.Cells.Item(24, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(24, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(24, 19).Value = 1250
$ws.Cells.Item(24, 20).Value = 16
